$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JAN-22")

# Row 4: new entry "2" dated 44565 (2022-01-04), RPA RLOGIC, GL summary task, 100%, Completed
# Clone number-format styles from existing analogous cells (B2 = date, E2 = percent)
# via copy/paste-formats so the shared cellXf entries are reused instead of duplicated.
$ws.Cells.Item(2, 2).Copy()
$ws.Cells.Item(4, 2).PasteSpecial(-4122)
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item(4, 5).PasteSpecial(-4122)
$ws.Cells.Item(3, 5).Copy()
$ws.Cells.Item(5, 5).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 44565
$ws.Cells.Item(4, 3).Value = "RPA RLOGIC"
$ws.Cells.Item(4, 4).Value = "1. Creating Borders and applying styles to the Summary of the General Ledger and it has been completed, implemented and tested"
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = "Completed"

# Row 5: continuation - Profit and Loss report task, 100%, Completed
$ws.Cells.Item(5, 4).Value = "2. Creating Borders and applying styles to the Summary of the Profit and Loss report and it has been completed, implemented and tested"
$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = "Completed"

# Update selection to D18
$ws.Range("D18").Select()
